# Updates the LR-pair sheet with new TPM-derived values and adds the
# 'Resolving-Mac' cluster, expanding the 3x3 sending/target grid to 4x4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T the NATMI numeric metrics (unchanged header row).
$data = @(
    @("ECs", "Tgfb2", "Tgfbr3", "ECs", 3, 1, 2.564153333333334, 7.692460000000001, 0.05249149542937438, 0.05249149542937438, 3, 1, 28.689524, 86.06857199999999, 0.2394085694101769, 0.2394085694101769, 73.56433859634666, 662.07904736712, 0.01256691382694736, 0.01256691382694736),
    @("ECs", "Tgfb2", "Tgfbr3", "FAPs", 3, 1, 2.564153333333334, 7.692460000000001, 0.05249149542937438, 0.05249149542937438, 3, 1, 62.46631733333334, 187.398952, 0.5212694246546397, 0.5212694246546395, 160.1732158113245, 1441.55894230192, 0.02736221162173163, 0.02736221162173163),
    @("ECs", "Tgfb2", "Tgfbr3", "MuSCs", 3, 1, 2.564153333333334, 7.692460000000001, 0.05249149542937438, 0.05249149542937438, 3, 1, 28.525746, 85.577238, 0.2380418734454457, 0.2380418734454457, 73.14438669172, 658.29948022548, 0.01249517391196133, 0.01249517391196133),
    @("ECs", "Tgfb2", "Tgfbr3", "Resolving-Mac", 3, 1, 2.564153333333334, 7.692460000000001, 0.05249149542937438, 0.05249149542937438, 3, 1, 0.1534046666666667, 0.460214, 0.001280132489737778, 0.001280132489737778, 0.3933530873822222, 3.54017778644, 0.0000671960687340642, 0.0000671960687340642),
    @("FAPs", "Tgfb2", "Tgfbr3", "ECs", 3, 1, 16.59481266666667, 49.78443799999999, 0.3397170215679993, 0.3397170215679993, 3, 1, 28.689524, 86.06857199999999, 0.2394085694101769, 0.2394085694101769, 476.0972762758373, 4284.875486482535, 0.08133116613788094, 0.08133116613788093),
    @("FAPs", "Tgfb2", "Tgfbr3", "FAPs", 3, 1, 16.59481266666667, 49.78443799999999, 0.3397170215679993, 0.3397170215679993, 3, 1, 62.46631733333334, 187.398952, 0.5212694246546397, 0.5212694246546395, 1036.61683412322, 9329.551507108976, 0.1770840963781388, 0.1770840963781387),
    @("FAPs", "Tgfb2", "Tgfbr3", "MuSCs", 3, 1, 16.59481266666667, 49.78443799999999, 0.3397170215679993, 0.3397170215679993, 3, 1, 28.525746, 85.577238, 0.2380418734454457, 0.2380418734454457, 473.3794110469159, 4260.414699422243, 0.08086687625535345, 0.08086687625535342),
    @("FAPs", "Tgfb2", "Tgfbr3", "Resolving-Mac", 3, 1, 16.59481266666667, 49.78443799999999, 0.3397170215679993, 0.3397170215679993, 3, 1, 0.1534046666666667, 0.460214, 0.001280132489737778, 0.001280132489737778, 2.545721705525778, 22.911495349732, 0.0004348827966261453, 0.0004348827966261452),
    @("MuSCs", "Tgfb2", "Tgfbr3", "ECs", 3, 1, 29.664466, 88.993398, 0.6072695268303631, 0.6072695268303631, 3, 1, 28.689524, 86.06857199999999, 0.2394085694101769, 0.2394085694101769, 851.0594092541838, 7659.534683287655, 0.1453855286648523, 0.1453855286648523),
    @("MuSCs", "Tgfb2", "Tgfbr3", "FAPs", 3, 1, 29.664466, 88.993398, 0.6072695268303631, 0.6072695268303631, 3, 1, 62.46631733333334, 187.398952, 0.5212694246546397, 0.5212694246546395, 1853.029946679877, 16677.2695201189, 0.3165510368611586, 0.3165510368611585),
    @("MuSCs", "Tgfb2", "Tgfbr3", "MuSCs", 3, 1, 29.664466, 88.993398, 0.6072695268303631, 0.6072695268303631, 3, 1, 28.525746, 85.577238, 0.2380418734454457, 0.2380418734454457, 846.201022341636, 7615.809201074723, 0.144555575853029, 0.144555575853029),
    @("MuSCs", "Tgfb2", "Tgfbr3", "Resolving-Mac", 3, 1, 29.664466, 88.993398, 0.6072695268303631, 0.6072695268303631, 3, 1, 0.1534046666666667, 0.460214, 0.001280132489737778, 0.001280132489737778, 4.550667518574667, 40.956007667172, 0.0007773854513232349, 0.0007773854513232349),
    @("Resolving-Mac", "Tgfb2", "Tgfbr3", "ECs", 1, 0.3333333333333333, 0.025497, 0.076491, 0.000521956172263265, 0.000521956172263265, 3, 1, 28.689524, 86.06857199999999, 0.2394085694101769, 0.2394085694101769, 0.7314967934279999, 6.583471140852, 0.0001249607804963602, 0.0001249607804963602),
    @("Resolving-Mac", "Tgfb2", "Tgfbr3", "FAPs", 1, 0.3333333333333333, 0.025497, 0.076491, 0.000521956172263265, 0.000521956172263265, 3, 1, 62.46631733333334, 187.398952, 0.5212694246546397, 0.5212694246546395, 1.592703693048, 14.334333237432, 0.0002720797936106102, 0.0002720797936106101),
    @("Resolving-Mac", "Tgfb2", "Tgfbr3", "MuSCs", 1, 0.3333333333333333, 0.025497, 0.076491, 0.000521956172263265, 0.000521956172263265, 3, 1, 28.525746, 85.577238, 0.2380418734454457, 0.2380418734454457, 0.727320945762, 6.545888511858, 0.0001242474251019614, 0.0001242474251019614),
    @("Resolving-Mac", "Tgfb2", "Tgfbr3", "Resolving-Mac", 1, 0.3333333333333333, 0.025497, 0.076491, 0.000521956172263265, 0.000521956172263265, 3, 1, 0.1534046666666667, 0.460214, 0.001280132489737778, 0.001280132489737778, 0.003911358786000001, 0.035202229074, 0.0000006681730543333738, 0.0000006681730543333738)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $row[$c]
    }
}

